$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell value updates per the upstream cryptos data refresh.
# Column D (Price) values that look like plain numbers must be forced to
# text (matching the original inlineStr storage) so Excel does not coerce
# them into numeric values, while columns B, C and E are already
# unambiguous text (URLs / names / percentages with "%" and padding spaces).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '21.736.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.29%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.540.62'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.88%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("E5").Value = '  -0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '290.22'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3901'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.69%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3180'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.70%  '
$ws.Range("E9").Value = '  +5.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07192'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.059'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.77%  '
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.635'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.64'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.16%  '
$ws.Range("B15").Value = 'WrappedEther'
$ws.Range("C15").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.574.73'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.58%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.629'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.86%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001102'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06583'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.99%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '83.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.16%  '
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("E21").Value = '  -4.31%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.396'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '21.731.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.363'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '146.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.54%  '
$ws.Range("E28").Value = '  -2.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.838'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.744.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '117.46'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.903'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.34%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9660'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -13.52%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08193'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.818'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06090'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.126'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02200'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2037'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.61%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.442'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -12.89%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.182'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '10.65'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5731'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -3.51%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.07'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.26%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.740'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5491'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.157'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '116.46'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.871'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06703'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.12%  '
